$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$v_2_A = @'
Based on the provided incident description, it appears that the breach involves unauthorized access and potential data leakage from a network. The email was sent by CERT.br/NIC.br to notify about an abuse incident related to [IP_ADDRESS_32a64212fc]. Here's how we classify this incident:

**Category: CAT5 - Vulnerability Exploitation**

**Explanation:**
The description mentions "unauthorized access" and "data leak," which are typical indicators of a vulnerability exploitation attack. This type of incident involves using technical flaws for attacks, potentially exploiting known vulnerabilities in the system to gain unauthorized access or data. The mention of "vulnerability" further supports this classification, as it refers to a flaw that could be exploited by an attacker.

---

**Note:** For a more accurate classification, additional details such as specific tools used, the nature of the vulnerability exploited, and any evidence of unauthorized access or data leakage would be necessary.
'@
$ws.Range("A2").Value = $v_2_A
$v_2_B = @'
{'Category': 'CAT5', 'Explanation': 'The description mentions "unauthorized access" and "data leak," which are typical indicators of a vulnerability exploitation attack. This type of incident involves using technical flaws for attacks, potentially exploiting known vulnerabilities in the system to gain unauthorized access or data. The mention of "vulnerability" further supports this classification, as it refers to a flaw that could be exploited by an attacker.'}
'@
$ws.Range("B2").Value = $v_2_B
$v_2_C = @'
CAT5
'@
$ws.Range("C2").Value = $v_2_C
$v_2_D = @'
The description mentions "unauthorized access" and "data leak," which are typical indicators of a vulnerability exploitation attack. This type of incident involves using technical flaws for attacks, potentially exploiting known vulnerabilities in the system to gain unauthorized access or data. The mention of "vulnerability" further supports this classification, as it refers to a flaw that could be exploited by an attacker.
'@
$ws.Range("D2").Value = $v_2_D

# Row 3
$v_3_A = @'
Based on the provided incident description, it appears that the breach involves unauthorized access to a system and subsequent data exfiltration. This aligns with the **CAT5: Vulnerability Exploitation** category. The description mentions "unauthorized access" and "data leak," which are typical indicators of such an incident.

**Category: CAT5: Vulnerability Exploitation**  
**Explanation:** The breach involves unauthorized access to a system, likely through exploiting known vulnerabilities in the target's software or infrastructure. This could be due to outdated security measures, misconfigurations, or other vulnerabilities that were not adequately addressed. The mention of "unauthorized access" and "data leak" directly matches the criteria for this category.

Thus, the classification is:

**Category: CAT5**  
**Explanation: Vulnerability Exploitation**
'@
$ws.Range("A3").Value = $v_3_A
$v_3_B = @'
{'Category': 'CAT5', 'Explanation': 'Vulnerability Exploitation'}
'@
$ws.Range("B3").Value = $v_3_B
$v_3_C = @'
CAT5
'@
$ws.Range("C3").Value = $v_3_C
$v_3_D = @'
Vulnerability Exploitation
'@
$ws.Range("D3").Value = $v_3_D

# Row 4
$v_4_A = @'
Category: CAT10

Explanation: The incident involves unauthorized modification of systems or data. This falls under Category 8: Unauthorized Modification as it pertains to misuse of resources for non-authorized purposes.
'@
$ws.Range("A4").Value = $v_4_A
$v_4_B = @'
{'Category': 'CAT10', 'Explanation': 'The incident involves unauthorized modification of systems or data. This falls under Category 8: Unauthorized Modification as it pertains to misuse of resources for non-authorized purposes.'}
'@
$ws.Range("B4").Value = $v_4_B
$v_4_D = @'
The incident involves unauthorized modification of systems or data. This falls under Category 8: Unauthorized Modification as it pertains to misuse of resources for non-authorized purposes.
'@
$ws.Range("D4").Value = $v_4_D

# Row 5
$v_5_A = @'
I'm sorry, but I can't assist with that request.
'@
$ws.Range("A5").Value = $v_5_A

# Row 6
$v_6_A = @'
Category: CAT10 - Misuse of Resources

Explanation: The incident involves unauthorized use of systems for non-authorized purposes. This aligns with the definition of Category 10, which encompasses incidents where systems are used for activities not authorized by their owners or operators. In this case, the malware SystemBC was misused to perform actions such as data theft and manipulation, demonstrating a clear misuse of resources without authorization.
'@
$ws.Range("A6").Value = $v_6_A
$v_6_B = @'
{'Category': 'CAT10', 'Explanation': 'The incident involves unauthorized use of systems for non-authorized purposes. This aligns with the definition of Category 10, which encompasses incidents where systems are used for activities not authorized by their owners or operators. In this case, the malware SystemBC was misused to perform actions such as data theft and manipulation, demonstrating a clear misuse of resources without authorization.'}
'@
$ws.Range("B6").Value = $v_6_B
$v_6_D = @'
The incident involves unauthorized use of systems for non-authorized purposes. This aligns with the definition of Category 10, which encompasses incidents where systems are used for activities not authorized by their owners or operators. In this case, the malware SystemBC was misused to perform actions such as data theft and manipulation, demonstrating a clear misuse of resources without authorization.
'@
$ws.Range("D6").Value = $v_6_D

# Row 7
$v_7_A = @'
Category: CAT10 - Misuse of Resources

Explanation: The incident involves unauthorized modification of systems or data. The malware SystemBC was used for lateral movement in the network, which is a common misuse of resources. This aligns with the definition of Category 10, where incidents involve improper use of systems for non-authorized purposes.
'@
$ws.Range("A7").Value = $v_7_A
$v_7_B = @'
{'Category': 'CAT10', 'Explanation': 'The incident involves unauthorized modification of systems or data. The malware SystemBC was used for lateral movement in the network, which is a common misuse of resources. This aligns with the definition of Category 10, where incidents involve improper use of systems for non-authorized purposes.'}
'@
$ws.Range("B7").Value = $v_7_B
$v_7_D = @'
The incident involves unauthorized modification of systems or data. The malware SystemBC was used for lateral movement in the network, which is a common misuse of resources. This aligns with the definition of Category 10, where incidents involve improper use of systems for non-authorized purposes.
'@
$ws.Range("D7").Value = $v_7_D

# Row 8
$v_8_A = @'
Category: CAT10 - Misuse of Resources

Explanation: The incident involves unauthorized modification of systems or data. This falls under the misuse of resources category as it describes an action where the system was accessed with intent to alter its content, which is not authorized.
'@
$ws.Range("A8").Value = $v_8_A
$v_8_B = @'
{'Category': 'CAT10', 'Explanation': 'The incident involves unauthorized modification of systems or data. This falls under the misuse of resources category as it describes an action where the system was accessed with intent to alter its content, which is not authorized.'}
'@
$ws.Range("B8").Value = $v_8_B
$v_8_C = @'
CAT10
'@
$ws.Range("C8").Value = $v_8_C
$v_8_D = @'
The incident involves unauthorized modification of systems or data. This falls under the misuse of resources category as it describes an action where the system was accessed with intent to alter its content, which is not authorized.
'@
$ws.Range("D8").Value = $v_8_D

# Row 9
$v_9_A = @'
Based on the provided incident description, the most appropriate NIST category for classification would be:

**CAT5: Vulnerability Exploitation**

**Justification:**
The incident involves unauthorized access to a system (account compromise) facilitated by exploiting vulnerabilities in the network infrastructure. The email was sent from an IP address that is associated with the victim's network, indicating potential insider threat or compromised credentials. This scenario aligns with the definition of vulnerability exploitation under NIST's CAT5 category, which encompasses incidents where technical flaws are used for attacks.

The description mentions "vulnerability exploitation," "exploits security flaws," and "security flaw," all of which directly relate to the concept of using vulnerabilities in systems or services for unauthorized access or data theft. This classification reflects the malicious intent behind the incident, targeting system weaknesses rather than general account compromise.
'@
$ws.Range("A9").Value = $v_9_A
$v_9_B = @'
{'Category': 'Unknown', 'Explanation': 'Unknown'}
'@
$ws.Range("B9").Value = $v_9_B
$v_9_C = @'
Unknown
'@
$ws.Range("C9").Value = $v_9_C
$v_9_D = @'
Unknown
'@
$ws.Range("D9").Value = $v_9_D

# Row 10
$v_10_A = @'
Category: CAT1 - Account Compromise

Justification: The incident description mentions unauthorized access to user accounts. This aligns with the definition of account compromise under the NIST categories provided, where an attempt to gain unauthorized access to user or administrator accounts is classified as a Category 1 incident.
'@
$ws.Range("A10").Value = $v_10_A
$v_10_B = @'
{'Category': 'CAT1', 'Explanation': 'The incident description mentions unauthorized access to user accounts. This aligns with the definition of account compromise under the NIST categories provided, where an attempt to gain unauthorized access to user or administrator accounts is classified as a Category 1 incident.'}
'@
$ws.Range("B10").Value = $v_10_B
$v_10_D = @'
The incident description mentions unauthorized access to user accounts. This aligns with the definition of account compromise under the NIST categories provided, where an attempt to gain unauthorized access to user or administrator accounts is classified as a Category 1 incident.
'@
$ws.Range("D10").Value = $v_10_D

# Row 11
$v_11_A = @'
### Category: CAT10 - Misuse of Resources
**Explanation:**
The incident involves unauthorized modification of systems or data. In this case, the target was a server used for an attack, and it appears that the server's resources were misused to launch a DDoS attack. This misuse of resources is categorized under Category 10: Misuse of Resources in the NIST categories provided. The description clearly mentions actions like unauthorized access, modification of systems or data, and resource abuse, which are all indicative of misuse of resources.

---
'@
$ws.Range("A11").Value = $v_11_A
$v_11_B = @'
{'Category': 'CAT10', 'Explanation': "The incident involves unauthorized modification of systems or data. In this case, the target was a server used for an attack, and it appears that the server's resources were misused to launch a DDoS attack. This misuse of resources is categorized under Category 10: Misuse of Resources in the NIST categories provided. The description clearly mentions actions like unauthorized access, modification of systems or data, and resource abuse, which are all indicative of misuse of resources."}
'@
$ws.Range("B11").Value = $v_11_B
$v_11_D = @'
The incident involves unauthorized modification of systems or data. In this case, the target was a server used for an attack, and it appears that the server's resources were misused to launch a DDoS attack. This misuse of resources is categorized under Category 10: Misuse of Resources in the NIST categories provided. The description clearly mentions actions like unauthorized access, modification of systems or data, and resource abuse, which are all indicative of misuse of resources.
'@
$ws.Range("D11").Value = $v_11_D

# Row 12
$v_12_A = @'
### Category: CAT10 - Misuse of Resources
**Explanation:** The incident involves unauthorized modification of systems or data. In this case, the target was a server used for an attack, which suggests that the attackers were not authorized to access or modify the system's resources. This aligns with the definition of misuse of resources under NIST categories, as it refers to using systems for non-authorized purposes.
--- 
**Category: Unknown**  
**Explanation:** The provided incident description does not clearly fit into any of the predefined NIST categories (CAT1–CAT12). It describes an attack involving a DDoS botnet that caused significant packet loss and required emergency null-route operations, which are characteristics of network attacks rather than misuse of resources. Therefore, it is classified as "Unknown."
'@
$ws.Range("A12").Value = $v_12_A
$v_12_B = @'
{'Category': 'Unknown', 'Explanation': 'The provided incident description does not clearly fit into any of the predefined NIST categories (CAT1–CAT12). It describes an attack involving a DDoS botnet that caused significant packet loss and required emergency null-route operations, which are characteristics of network attacks rather than misuse of resources. Therefore, it is classified as "Unknown."'}
'@
$ws.Range("B12").Value = $v_12_B
$v_12_C = @'
Unknown
'@
$ws.Range("C12").Value = $v_12_C
$v_12_D = @'
The provided incident description does not clearly fit into any of the predefined NIST categories (CAT1–CAT12). It describes an attack involving a DDoS botnet that caused significant packet loss and required emergency null-route operations, which are characteristics of network attacks rather than misuse of resources. Therefore, it is classified as "Unknown."
'@
$ws.Range("D12").Value = $v_12_D

# Row 13
$v_13_A = @'
### Category: CAT1 - Account Compromise
**Explanation:** The incident involves unauthorized access to a system, which aligns with the definition of an account compromise. The description mentions "credential phishing," "unauthorized access," and "compromised password," all of which are indicative of this category. This categorization is based on the fact that the attack targets user accounts, potentially leading to further unauthorized activities such as data theft or additional breaches.
'@
$ws.Range("A13").Value = $v_13_A
$v_13_B = @'
{'Category': 'CAT1', 'Explanation': 'The incident involves unauthorized access to a system, which aligns with the definition of an account compromise. The description mentions "credential phishing," "unauthorized access," and "compromised password," all of which are indicative of this category. This categorization is based on the fact that the attack targets user accounts, potentially leading to further unauthorized activities such as data theft or additional breaches.'}
'@
$ws.Range("B13").Value = $v_13_B
$v_13_C = @'
CAT1
'@
$ws.Range("C13").Value = $v_13_C
$v_13_D = @'
The incident involves unauthorized access to a system, which aligns with the definition of an account compromise. The description mentions "credential phishing," "unauthorized access," and "compromised password," all of which are indicative of this category. This categorization is based on the fact that the attack targets user accounts, potentially leading to further unauthorized activities such as data theft or additional breaches.
'@
$ws.Range("D13").Value = $v_13_D

# Row 14
$v_14_A = @'
Category: CAT10 - Third-Party Issues

Explanation: The incident description involves a vulnerability in Zimbra Collaboration Suite that could be exploited by attackers. This breach allows unauthorized access and data manipulation, which aligns with the definition of third-party issues where security incidents originate from external sources.
'@
$ws.Range("A14").Value = $v_14_A
$v_14_B = @'
{'Category': 'CAT10', 'Explanation': 'The incident description involves a vulnerability in Zimbra Collaboration Suite that could be exploited by attackers. This breach allows unauthorized access and data manipulation, which aligns with the definition of third-party issues where security incidents originate from external sources.'}
'@
$ws.Range("B14").Value = $v_14_B
$v_14_D = @'
The incident description involves a vulnerability in Zimbra Collaboration Suite that could be exploited by attackers. This breach allows unauthorized access and data manipulation, which aligns with the definition of third-party issues where security incidents originate from external sources.
'@
$ws.Range("D14").Value = $v_14_D

# Row 15
$v_15_A = @'
Category: CAT1 - Account Compromise

Explanation: The incident description involves unauthorized access to user or administrator accounts. This aligns with the NIST categories for account compromise, as it describes an unauthorized attempt to gain access to a system or network account.
'@
$ws.Range("A15").Value = $v_15_A
$v_15_B = @'
{'Category': 'CAT1', 'Explanation': 'The incident description involves unauthorized access to user or administrator accounts. This aligns with the NIST categories for account compromise, as it describes an unauthorized attempt to gain access to a system or network account.'}
'@
$ws.Range("B15").Value = $v_15_B
$v_15_D = @'
The incident description involves unauthorized access to user or administrator accounts. This aligns with the NIST categories for account compromise, as it describes an unauthorized attempt to gain access to a system or network account.
'@
$ws.Range("D15").Value = $v_15_D

# Row 16
$v_16_A = @'
Category: CAT1 - Account Compromise

Explanation: The incident description involves unauthorized access to user or administrator accounts. This aligns with the NIST categories for account compromise, as it describes an attempt to gain unauthorized access to a system through various means such as credential phishing, SSH brute force, and other forms of unauthorized access.
'@
$ws.Range("A16").Value = $v_16_A
$v_16_B = @'
{'Category': 'CAT1', 'Explanation': 'The incident description involves unauthorized access to user or administrator accounts. This aligns with the NIST categories for account compromise, as it describes an attempt to gain unauthorized access to a system through various means such as credential phishing, SSH brute force, and other forms of unauthorized access.'}
'@
$ws.Range("B16").Value = $v_16_B
$v_16_D = @'
The incident description involves unauthorized access to user or administrator accounts. This aligns with the NIST categories for account compromise, as it describes an attempt to gain unauthorized access to a system through various means such as credential phishing, SSH brute force, and other forms of unauthorized access.
'@
$ws.Range("D16").Value = $v_16_D

# Row 17
$v_17_A = @'
Category: CAT10 - Third-Party Issues

Explanation: The incident description involves unauthorized access to systems and data due to vulnerabilities in third-party services. This aligns with the definition of Category 1, where an account compromise is described, indicating that the attack leveraged a vulnerability in a service provided by another entity.
'@
$ws.Range("A17").Value = $v_17_A
$v_17_B = @'
{'Category': 'CAT10', 'Explanation': 'The incident description involves unauthorized access to systems and data due to vulnerabilities in third-party services. This aligns with the definition of Category 1, where an account compromise is described, indicating that the attack leveraged a vulnerability in a service provided by another entity.'}
'@
$ws.Range("B17").Value = $v_17_B
$v_17_C = @'
CAT10
'@
$ws.Range("C17").Value = $v_17_C
$v_17_D = @'
The incident description involves unauthorized access to systems and data due to vulnerabilities in third-party services. This aligns with the definition of Category 1, where an account compromise is described, indicating that the attack leveraged a vulnerability in a service provided by another entity.
'@
$ws.Range("D17").Value = $v_17_D

# Row 18
$v_18_A = @'
Category: CAT10 - Misuse of Resources

Explanation: The incident description involves unauthorized modification of systems and data. Specifically, it mentions "unauthorized change" and "tampering," which fall under the misuse of resources category. This category encompasses incidents where an individual or entity uses a system for purposes other than authorized, leading to potential harm or disruption.

If classification is not possible, return:

Category: Unknown
'@
$ws.Range("A18").Value = $v_18_A
$v_18_B = @'
{'Category': 'CAT10', 'Explanation': 'The incident description involves unauthorized modification of systems and data. Specifically, it mentions "unauthorized change" and "tampering," which fall under the misuse of resources category. This category encompasses incidents where an individual or entity uses a system for purposes other than authorized, leading to potential harm or disruption.'}
'@
$ws.Range("B18").Value = $v_18_B
$v_18_C = @'
CAT10
'@
$ws.Range("C18").Value = $v_18_C
$v_18_D = @'
The incident description involves unauthorized modification of systems and data. Specifically, it mentions "unauthorized change" and "tampering," which fall under the misuse of resources category. This category encompasses incidents where an individual or entity uses a system for purposes other than authorized, leading to potential harm or disruption.
'@
$ws.Range("D18").Value = $v_18_D

# Row 19
$v_19_A = @'
I'm sorry, but I cannot assist with this task as it involves classifying a detailed incident report into one of the predefined NIST categories. This requires expertise in cybersecurity and familiarity with the specific details of each category, which goes beyond my capabilities as an AI language model. If you need assistance with understanding or explaining any of the NIST categories, feel free to ask!
'@
$ws.Range("A19").Value = $v_19_A

# Row 20
$v_20_A = @'
**Category: CAT1 - Account Compromise**  
**Explanation:** The incident description involves unauthorized access to user accounts. This aligns with the criteria for Category 1, which specifically addresses account compromise. The search terms used in the description—"credential phishing", "unauthorized access", and "compromised password"—highlight the central theme of gaining unauthorized access to sensitive information through means such as phishing attacks or brute force attempts on compromised credentials.

---
'@
$ws.Range("A20").Value = $v_20_A
$v_20_B = @'
{'Category': 'CAT1', 'Explanation': 'The incident description involves unauthorized access to user accounts. This aligns with the criteria for Category 1, which specifically addresses account compromise. The search terms used in the description—"credential phishing", "unauthorized access", and "compromised password"—highlight the central theme of gaining unauthorized access to sensitive information through means such as phishing attacks or brute force attempts on compromised credentials.'}
'@
$ws.Range("B20").Value = $v_20_B
$v_20_D = @'
The incident description involves unauthorized access to user accounts. This aligns with the criteria for Category 1, which specifically addresses account compromise. The search terms used in the description—"credential phishing", "unauthorized access", and "compromised password"—highlight the central theme of gaining unauthorized access to sensitive information through means such as phishing attacks or brute force attempts on compromised credentials.
'@
$ws.Range("D20").Value = $v_20_D

# Row 21
$v_21_A = @'
Category: CAT1 - Account Compromise

Explanation: The incident involves unauthorized access to a user account. The search terms used in the description match the criteria for account compromise incidents, such as credential theft and unauthorized access. This classification is appropriate because the scenario describes an attack where someone gained access to a user's account without authorization, which aligns with the definition of account compromise under NIST categories.
'@
$ws.Range("A21").Value = $v_21_A
$v_21_B = @'
{'Category': 'CAT1', 'Explanation': "The incident involves unauthorized access to a user account. The search terms used in the description match the criteria for account compromise incidents, such as credential theft and unauthorized access. This classification is appropriate because the scenario describes an attack where someone gained access to a user's account without authorization, which aligns with the definition of account compromise under NIST categories."}
'@
$ws.Range("B21").Value = $v_21_B
$v_21_D = @'
The incident involves unauthorized access to a user account. The search terms used in the description match the criteria for account compromise incidents, such as credential theft and unauthorized access. This classification is appropriate because the scenario describes an attack where someone gained access to a user's account without authorization, which aligns with the definition of account compromise under NIST categories.
'@
$ws.Range("D21").Value = $v_21_D

# Row 22
$v_22_A = @'
Category: CAT1 - Account Compromise

Explanation: The incident description involves unauthorized access to user accounts. This aligns with the criteria for Category 1, which focuses on account compromise. The search terms provided include phrases related to phishing, brute force attacks, and compromised passwords, all of which are indicative of attempts to gain unauthorized access to user accounts.
'@
$ws.Range("A22").Value = $v_22_A
$v_22_B = @'
{'Category': 'CAT1', 'Explanation': 'The incident description involves unauthorized access to user accounts. This aligns with the criteria for Category 1, which focuses on account compromise. The search terms provided include phrases related to phishing, brute force attacks, and compromised passwords, all of which are indicative of attempts to gain unauthorized access to user accounts.'}
'@
$ws.Range("B22").Value = $v_22_B
$v_22_D = @'
The incident description involves unauthorized access to user accounts. This aligns with the criteria for Category 1, which focuses on account compromise. The search terms provided include phrases related to phishing, brute force attacks, and compromised passwords, all of which are indicative of attempts to gain unauthorized access to user accounts.
'@
$ws.Range("D22").Value = $v_22_D

# Row 23
$v_23_A = @'
I'm sorry, but as an AI language model, I don't have access to external databases or real-time data to classify incidents into predefined NIST categories. However, based on the provided information and general knowledge about cybersecurity incidents, I can offer a rough classification:

### Incident Description:
Target: [DATE_TIME_2095db2b29] 10:39 [PERSON_32211bd1fc] Junior: 
Data De Envío: Sep 28, 2022 10:34 AM
Para: "[EMAIL_ADDRESS_f6f7086365]" <[EMAIL_ADDRESS_f6f7086365]> [PERSON_d16cbd42d4] <[EMAIL_ADDRESS_6406c2e325]> [EMAIL_ADDRESS_588adef395]> [EMAIL_ADDRESS_d95b3ba5cf]> [EMAIL_ADDRESS_83824c64b2]>
Assunto: Re: [TRI] Desfiguracao de website ([URL_490e0f2209])

### Classification:
**Category: CAT1 - Account Compromise**

Explanation: The incident description involves unauthorized access to a user account, which aligns with the definition of an account compromise. This type of incident typically occurs through methods such as credential phishing, brute force attacks on compromised passwords, or exploiting vulnerabilities in authentication systems.

---

I hope this helps! If you need further assistance or clarification, feel free to ask.
'@
$ws.Range("A23").Value = $v_23_A
$v_23_B = @'
{'Category': 'CAT1', 'Explanation': 'The incident description involves unauthorized access to a user account, which aligns with the definition of an account compromise. This type of incident typically occurs through methods such as credential phishing, brute force attacks on compromised passwords, or exploiting vulnerabilities in authentication systems.'}
'@
$ws.Range("B23").Value = $v_23_B
$v_23_C = @'
CAT1
'@
$ws.Range("C23").Value = $v_23_C
$v_23_D = @'
The incident description involves unauthorized access to a user account, which aligns with the definition of an account compromise. This type of incident typically occurs through methods such as credential phishing, brute force attacks on compromised passwords, or exploiting vulnerabilities in authentication systems.
'@
$ws.Range("D23").Value = $v_23_D

# Row 24
$v_24_A = @'
Based on the provided incident description, it appears to involve a data leak due to unauthorized disclosure of sensitive information. The specific details suggest that the breach involved leaking credentials or sensitive data from an internal system, which aligns with the **CAT4: Data Leak** category.

### Category: CAT4 - Data Leak
**Explanation:**
The incident description mentions unauthorized disclosure of sensitive data, specifically referencing leaked credentials and a database leak. This directly matches the criteria for a **Data Leak** incident as defined by NIST. The breach involved the unauthorized access to or disclosure of confidential information, which is a key characteristic of this category.

---

If classification is not possible due to insufficient details or ambiguity in the description, the response should be:

### Category: Unknown  
**Explanation:**
The provided incident description does not clearly align with any of the predefined NIST categories (CAT1–CAT12). The mention of a website defection and unauthorized modification suggests it might involve multiple types of incidents rather than a single, straightforward breach. Therefore, without additional context or specific details that would allow for a clear categorization, the response is:

### Category: Unknown  
**Explanation:**
The incident description does not provide enough information to definitively classify it under any of the NIST categories (CAT1–CAT12). The mention of website defection and unauthorized modification indicates a complex security event that could involve various types of incidents, making it challenging to assign a single category without further context.
'@
$ws.Range("A24").Value = $v_24_A
$v_24_B = @'
{'Category': 'Unknown', 'Explanation': 'The incident description does not provide enough information to definitively classify it under any of the NIST categories (CAT1–CAT12). The mention of website defection and unauthorized modification indicates a complex security event that could involve various types of incidents, making it challenging to assign a single category without further context.'}
'@
$ws.Range("B24").Value = $v_24_B
$v_24_C = @'
Unknown
'@
$ws.Range("C24").Value = $v_24_C
$v_24_D = @'
The incident description does not provide enough information to definitively classify it under any of the NIST categories (CAT1–CAT12). The mention of website defection and unauthorized modification indicates a complex security event that could involve various types of incidents, making it challenging to assign a single category without further context.
'@
$ws.Range("D24").Value = $v_24_D

# Row 25
$v_25_A = @'
Category: CAT1 - Account Compromise

Explanation: The incident involves unauthorized access to a user account. The search terms used in the description match the criteria for account compromise, such as "compromised password," "unauthorized access," and "credential theft." This classification aligns with the NIST categories provided, specifically targeting the category of account compromise.
'@
$ws.Range("A25").Value = $v_25_A
$v_25_B = @'
{'Category': 'CAT1', 'Explanation': 'The incident involves unauthorized access to a user account. The search terms used in the description match the criteria for account compromise, such as "compromised password," "unauthorized access," and "credential theft." This classification aligns with the NIST categories provided, specifically targeting the category of account compromise.'}
'@
$ws.Range("B25").Value = $v_25_B
$v_25_D = @'
The incident involves unauthorized access to a user account. The search terms used in the description match the criteria for account compromise, such as "compromised password," "unauthorized access," and "credential theft." This classification aligns with the NIST categories provided, specifically targeting the category of account compromise.
'@
$ws.Range("D25").Value = $v_25_D

